$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StorePage")

# Existing columns (A/B already hold productName/sortingText, C/H already
# hold name/state text shared with other sheets) - fill column-by-column so
# newly introduced shared strings are interned in the same order as the
# target workbook.
$ws.Range("C1").Value = "name"
$ws.Range("C2").Value = "prashant"

$ws.Range("D1").Value = "lastName"
$ws.Range("D2").Value = "More"

$ws.Range("E1").Value = "country"
$ws.Range("E2").Value = "India"

$ws.Range("F1").Value = "address"
$ws.Range("F2").Value = "Gaoan Bhag"

$ws.Range("G1").Value = "city"
$ws.Range("G2").Value = "Kolhapur"

$ws.Range("H1").Value = "state"
$ws.Range("H2").Value = "Maharashtra"

$ws.Range("I1").Value = "email"
$ws.Range("I2").Value = "test@gmail.com"

# Hyperlink the email cell, matching the style used on the "NEW CUSTOMER" sheet.
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:test@gmail.com")
$ws.Range("I2").Style = "Hyperlink"

# Update the active selection to match the edited cell.
$ws.Range("I2").Select()
